$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for the 2022-Q4 quarter at the
#    top of the data (row 2), shifting the existing Q3/Q2/Q1 rows down, and
#    renumber the index column (A) so it stays 0,1,2,3.
# ---------------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")

$sheetTotal.Rows.Item(2).Insert()
$sheetTotal.Range("B2:D2").ClearFormats()

# Column A carries a bordered/centered style (index used by the existing
# index cells) - copy that formatting onto the freshly inserted A2 cell.
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("C2").Value = 23
$sheetTotal.Range("D2").Value = 1.6

$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("A4").Value = 2
$sheetTotal.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet right after "总计", ahead of the other
#    quarter sheets (2022-Q3 / 2022-Q2 / 2022-Q1), and populate it with the
#    quarter's fund holdings table.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $sheetTotal)
$newSheet.Name = "2022-Q4"

# Bring over the header / index-column formatting from an existing quarter
# sheet so the new sheet matches the workbook's established look.
$templateSheet = $wb.Worksheets.Item("2022-Q3")

$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A24").PasteSpecial(-4122)

    $data = New-Object 'object[,]' 24,8
    $data[0,1] = "基金代码"
    $data[0,2] = "基金名称"
    $data[0,3] = "基金规模"
    $data[0,4] = "股票总仓位"
    $data[0,5] = "仓位占比"
    $data[0,6] = "持有市值(亿元)"
    $data[0,7] = "仓位排名"
    $data[1,0] = 0
    $data[1,1] = "'050008"
    $data[1,2] = "博时第三产业成长混合"
    $data[1,3] = "'9.01"
    $data[1,4] = "'87.80"
    $data[1,5] = "'5.39"
    $data[1,6] = "'0.4856"
    $data[1,7] = 2
    $data[2,0] = 1
    $data[2,1] = "'011069"
    $data[2,2] = "工银成长精选混合A"
    $data[2,3] = "'12.99"
    $data[2,4] = "'91.27"
    $data[2,5] = "'2.98"
    $data[2,6] = "'0.3871"
    $data[2,7] = 9
    $data[3,0] = 2
    $data[3,1] = "'013417"
    $data[3,2] = "博时核心资产精选混合A"
    $data[3,3] = "'6.96"
    $data[3,4] = "'80.13"
    $data[3,5] = "'3.22"
    $data[3,6] = "'0.2241"
    $data[3,7] = 7
    $data[4,0] = 3
    $data[4,1] = "'012153"
    $data[4,2] = "博时研究慧选混合A"
    $data[4,3] = "'1.15"
    $data[4,4] = "'69.54"
    $data[4,5] = "'5.70"
    $data[4,6] = "'0.0656"
    $data[4,7] = 3
    $data[5,0] = 4
    $data[5,1] = "'003032"
    $data[5,2] = "平安医疗健康灵活配置混合"
    $data[5,3] = "'1.42"
    $data[5,4] = "'94.02"
    $data[5,5] = "'4.41"
    $data[5,6] = "'0.0626"
    $data[5,7] = 5
    $data[6,0] = 5
    $data[6,1] = "'160512"
    $data[6,2] = "博时卓越品牌混合（LOF）"
    $data[6,3] = "'1.84"
    $data[6,4] = "'68.12"
    $data[6,5] = "'3.30"
    $data[6,6] = "'0.0607"
    $data[6,7] = 7
    $data[7,0] = 6
    $data[7,1] = "'005265"
    $data[7,2] = "博时厚泽回报灵活配置混合A"
    $data[7,3] = "'1.63"
    $data[7,4] = "'77.88"
    $data[7,5] = "'3.44"
    $data[7,6] = "'0.0561"
    $data[7,7] = 6
    $data[8,0] = 7
    $data[8,1] = "'011070"
    $data[8,2] = "工银成长精选混合C"
    $data[8,3] = "'1.82"
    $data[8,4] = "'91.27"
    $data[8,5] = "'2.98"
    $data[8,6] = "'0.0542"
    $data[8,7] = 9
    $data[9,0] = 8
    $data[9,1] = "'015902"
    $data[9,2] = "博时优质精选混合A"
    $data[9,3] = "'3.15"
    $data[9,4] = "'30.80"
    $data[9,5] = "'1.69"
    $data[9,6] = "'0.0532"
    $data[9,7] = 5
    $data[10,0] = 9
    $data[10,1] = "'164826"
    $data[10,2] = "工银瑞信创业板两年定期开放混合A"
    $data[10,3] = "'1.64"
    $data[10,4] = "'80.79"
    $data[10,5] = "'2.88"
    $data[10,6] = "'0.0472"
    $data[10,7] = 10
    $data[11,0] = 10
    $data[11,1] = "'005266"
    $data[11,2] = "博时厚泽回报灵活配置混合C"
    $data[11,3] = "'0.90"
    $data[11,4] = "'77.88"
    $data[11,5] = "'3.44"
    $data[11,6] = "'0.0310"
    $data[11,7] = 6
    $data[12,0] = 11
    $data[12,1] = "'162212"
    $data[12,2] = "泰达宏利红利先锋混合A"
    $data[12,3] = "'0.59"
    $data[12,4] = "'93.38"
    $data[12,5] = "'3.59"
    $data[12,6] = "'0.0212"
    $data[12,7] = 10
    $data[13,0] = 12
    $data[13,1] = "'013418"
    $data[13,2] = "博时核心资产精选混合C"
    $data[13,3] = "'0.38"
    $data[13,4] = "'80.13"
    $data[13,5] = "'3.22"
    $data[13,6] = "'0.0122"
    $data[13,7] = 7
    $data[14,0] = 13
    $data[14,1] = "'012154"
    $data[14,2] = "博时研究慧选混合C"
    $data[14,3] = "'0.16"
    $data[14,4] = "'69.54"
    $data[14,5] = "'5.70"
    $data[14,6] = "'0.0091"
    $data[14,7] = 3
    $data[15,0] = 14
    $data[15,1] = "'002068"
    $data[15,2] = "东方多策略灵活配置混合C"
    $data[15,3] = "'0.26"
    $data[15,4] = "'87.87"
    $data[15,5] = "'2.85"
    $data[15,6] = "'0.0074"
    $data[15,7] = 6
    $data[16,0] = 15
    $data[16,1] = "'015903"
    $data[16,2] = "博时优质精选混合C"
    $data[16,3] = "'0.34"
    $data[16,4] = "'30.80"
    $data[16,5] = "'1.69"
    $data[16,6] = "'0.0057"
    $data[16,7] = 5
    $data[17,0] = 16
    $data[17,1] = "'002598"
    $data[17,2] = "平安消费精选混合A"
    $data[17,3] = "'0.13"
    $data[17,4] = "'93.51"
    $data[17,5] = "'4.39"
    $data[17,6] = "'0.0057"
    $data[17,7] = 5
    $data[18,0] = 17
    $data[18,1] = "'002599"
    $data[18,2] = "平安消费精选混合C"
    $data[18,3] = "'0.10"
    $data[18,4] = "'93.51"
    $data[18,5] = "'4.39"
    $data[18,6] = "'0.0044"
    $data[18,7] = 5
    $data[19,0] = 18
    $data[19,1] = "'010889"
    $data[19,2] = "工银瑞信创业板两年定期开放混合C"
    $data[19,3] = "'0.14"
    $data[19,4] = "'80.79"
    $data[19,5] = "'2.88"
    $data[19,6] = "'0.0040"
    $data[19,7] = 10
    $data[20,0] = 19
    $data[20,1] = "'005209"
    $data[20,2] = "东吴双三角股票A"
    $data[20,3] = "'0.09"
    $data[20,4] = "'92.11"
    $data[20,5] = "'2.97"
    $data[20,6] = "'0.0027"
    $data[20,7] = 10
    $data[21,0] = 20
    $data[21,1] = "'005210"
    $data[21,2] = "东吴双三角股票C"
    $data[21,3] = "'0.09"
    $data[21,4] = "'92.11"
    $data[21,5] = "'2.97"
    $data[21,6] = "'0.0027"
    $data[21,7] = 10
    $data[22,0] = 21
    $data[22,1] = "'400023"
    $data[22,2] = "东方多策略灵活配置混合A"
    $data[22,3] = "'0.03"
    $data[22,4] = "'87.87"
    $data[22,5] = "'2.85"
    $data[22,6] = "'0.0009"
    $data[22,7] = 6
    $data[23,0] = 22
    $data[23,1] = "'015619"
    $data[23,2] = "泰达宏利红利先锋混合C"
    $data[23,3] = "'0.00"
    $data[23,4] = "'93.38"
    $data[23,5] = "'3.59"
    $data[23,6] = 0
    $data[23,7] = 10

$newSheet.Range("A1:H24").Value = $data

Write-Output "Inserted 2022-Q4 sheet and updated summary sheet"
